# Apply "gh-pages output generated at 456a3b4" update.
# Workbook layout: 1=展览(Exhibitions) 2=演出(Performances) 3=本地生活(Local Life) 4=全部类型(All types)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 - bump "想去人数" (want-to-go count) in column F for the rows
# whose underlying event page count changed.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 7245
$ws1.Range("F4").Value  = 3501
$ws1.Range("F6").Value  = 3831
$ws1.Range("F7").Value  = 65
$ws1.Range("F8").Value  = 73
$ws1.Range("F10").Value = 98
$ws1.Range("F11").Value = 140
$ws1.Range("F14").Value = 127
$ws1.Range("F16").Value = 13
$ws1.Range("F19").Value = 4099
$ws1.Range("F24").Value = 1643
$ws1.Range("F27").Value = 3006
$ws1.Range("F28").Value = 2196
$ws1.Range("F31").Value = 91
$ws1.Range("F33").Value = 85
$ws1.Range("F36").Value = 4270
$ws1.Range("F37").Value = 463
$ws1.Range("F41").Value = 785
$ws1.Range("F42").Value = 197
$ws1.Range("F44").Value = 1621
$ws1.Range("F45").Value = 260
$ws1.Range("F46").Value = 29
$ws1.Range("F47").Value = 601
$ws1.Range("F48").Value = 712

# ---------------------------------------------------------------------------
# Sheet 2: 演出 - same kind of refresh.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value  = 438
$ws2.Range("F16").Value = 566

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 - refresh existing row, then append the new event row 3.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 164

# New row 3 data.
$ws3.Range("A3").Value = 2

# Column B holds a plain "yyyy-mm-dd" text label in this sheet (see B2), not a
# real date, so force text formatting via a leading quote to stop the COM
# layer from auto-converting the literal into a date serial, then strip the
# resulting cell style back off so the cell matches its plain, unstyled
# neighbours.
$ws3.Range("B3").Value = "'2024-06-16"
$ws3.Range("B3").Style = "Normal"

$ws3.Range("C3").Value = "北京·贰伊Lolita茶会"
$ws3.Range("D3").Value = "高碑店方家村甲西店记忆文创小镇D8号 格乐利雅·G婚礼艺术中心(朝阳店)"
$ws3.Range("E3").Value = "2024.06.16 14:00-06.16 17:30"
$ws3.Range("F3").Value = 0
$ws3.Range("G3").Value = 198
$ws3.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=86727"
$ws3.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202405/yZrFM4gf1717136810203.jpeg"

# Match the bold/bordered/centered look used by the other "#" column cells
# (copy format only from A2, which already carries that style).
$ws3.Range("A2").Copy()
$ws3.Range("A3").PasteSpecial(-4122)
$ws3.Range("A3").Value = 2

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 - aggregated view; mirrors the same refreshed counts.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 164
$ws4.Range("F5").Value  = 7245
$ws4.Range("F6").Value  = 3501
$ws4.Range("F7").Value  = 3501
$ws4.Range("F8").Value  = 3831
$ws4.Range("F9").Value  = 73
$ws4.Range("F11").Value = 98
$ws4.Range("F13").Value = 140
$ws4.Range("F16").Value = 127
$ws4.Range("F18").Value = 13
$ws4.Range("F21").Value = 4099
$ws4.Range("F28").Value = 1643
$ws4.Range("F31").Value = 3006
$ws4.Range("F32").Value = 2196
$ws4.Range("F35").Value = 91
$ws4.Range("F39").Value = 4270
$ws4.Range("F41").Value = 463
$ws4.Range("F45").Value = 785
$ws4.Range("F46").Value = 197
$ws4.Range("F47").Value = 1621
$ws4.Range("F48").Value = 260
$ws4.Range("F49").Value = 601
$ws4.Range("F50").Value = 712
